# Auto-generated cell updates for cryptos.xlsx "Price"/"Volume(1h)" refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.520.02"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "1.553.92"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").Value = "210.36"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("D8").Value = "24.29"
$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "1.776.13"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").Value = "1.554.28"
$ws.Range("E13").Value = "  -1.40%  "

$ws.Range("D14").Value = "28.492.34"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "61.21"
$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("D18").Value = "229.06"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").Value = "0.0₃0670"
$ws.Range("E20").Value = "  -2.28%  "

$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").Value = "3.89"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.71%  "

$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("D25").Value = "151.27"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").Value = "14.75"
$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("E29").Value = "  -2.20%  "

$ws.Range("D30").Value = "0.0458"
$ws.Range("E30").Value = "  -4.70%  "

$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").Value = "1.389.84"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("E34").Value = "  -2.73%  "

$ws.Range("E37").Value = "  -2.90%  "

$ws.Range("D38").Value = "2.64"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("E39").Value = "  -1.29%  "

$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "

$ws.Range("E43").Value = "  -1.11%  "

$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").Value = "64.24"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").Value = "5.29"

$ws.Range("D47").Value = "1.689.74"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("D48").Value = "0.869"
$ws.Range("E48").Value = "  -6.70%  "

$ws.Range("D49").Value = "43.52"
$ws.Range("E49").Value = "  +5.01%  "

$ws.Range("D50").Value = "85.16"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("E51").Value = "  +2.98%  "
